$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53..152 down to 54..153
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with its data
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C53").Value = 'Los Lagos'
$ws.Range("D53").Value = 44581
$ws.Range("E53").Value = 10
$ws.Range("F53").Value = 'Fruta'
$ws.Range("G53").Value = 100108
$ws.Range("H53").Value = 'Tropicales y subtropicales'
$ws.Range("I53").Value = 100108002
$ws.Range("J53").Value = 'Mango'
$ws.Range("K53").Value = 'Sin especificar'
$ws.Range("L53").Value = 'Primera'
$ws.Range("M53").Value = 120
$ws.Range("N53").Value = 8000
$ws.Range("O53").Value = 8500
$ws.Range("P53").Value = 8250
$ws.Range("Q53").Value = '$/bandeja 4 kilos'
$ws.Range("R53").Value = 'Perú'
$ws.Range("S53").Value = 2062
$ws.Range("T53").Value = 4

# Ensure date style for the new row's date cell matches the rest of column D
$ws.Range("D53").NumberFormat = $ws.Range("D54").NumberFormat
